# Daily attendance processing - 2026-01-01 16:37:38
# Normalizes the "Recorded By" column (G) ordering for specific known values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacement table: old "Recorded By" text -> new "Recorded By" text.
$replacements = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
